$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.623.66"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.343.58"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'543.58"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "'135.62"
$ws.Range("E6").Value = "  -5.94%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.523"
$ws.Range("E8").Value = "  -9.10%  "
$ws.Range("D9").Value = "2.342.78"
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").Value = "'24.46"
$ws.Range("E14").Value = "  -4.25%  "
$ws.Range("D15").Value = "2.766.40"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").Value = "60.510.17"
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("E17").Value = "  -2.97%  "
$ws.Range("D18").Value = "2.342.66"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "'317.56"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "'6.53"
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").Value = "'1.70"
$ws.Range("E25").Value = "  -4.82%  "
$ws.Range("E26").Value = "  +7.57%  "
$ws.Range("D27").Value = "'7.90"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "'496.23"
$ws.Range("E28").Value = "  -3.79%  "
$ws.Range("E29").Value = "  -5.56%  "
$ws.Range("D30").Value = "0.0₃0855"
$ws.Range("E30").Value = "  -10.04%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("E32").Value = "  -3.71%  "
$ws.Range("E33").Value = "  -4.06%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").Value = "'0.374"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("E37").Value = "  +2.37%  "
$ws.Range("E38").Value = "  -5.98%  "
$ws.Range("D39").Value = "'1.79"
$ws.Range("E39").Value = "  +4.27%  "
$ws.Range("D40").Value = "'141.13"
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "'141.27"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").Value = "'2.05"
$ws.Range("D46").Value = "'0.0509"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "'18.92"
$ws.Range("E47").Value = "  -8.44%  "
$ws.Range("E48").Value = "  -2.51%  "
$ws.Range("E49").Value = "  -2.89%  "
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("D51").Value = "'16.34"
$ws.Range("E51").Value = "  -2.42%  "
